# Updates the "cryptos" price/volume table (commit: "Updated cryptos list" GitHub Action).
# Columns: D = Price (text), E = Volume(1h) (text, "  +x.xx%  ").
#
# D holds plain text in the source sheet, and several new Price strings still look
# numeric (e.g. "0.9995", "146.00", "0.000008474"). Assigning those to Range.Value
# directly would make Excel coerce them into real numbers (dropping trailing zeros /
# using scientific notation), so for those we prefix with a leading apostrophe -
# the standard Excel "force text" entry method - which keeps the literal text.
# Values that are unambiguously non-numeric (e.g. "27.070.44", multiple dots) are
# assigned as-is since Excel already stores them as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.070.44'
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").Value = '1.890.69'
$ws.Range("E3").Value = '  +1.64%  '
$ws.Range("D4").Value = '''0.9995'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''306.96'
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("D6").Value = '''0.9994'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '''0.5142'
$ws.Range("E7").Value = '  +1.90%  '
$ws.Range("D8").Value = '''0.3754'
$ws.Range("E8").Value = '  +3.21%  '
$ws.Range("D9").Value = '''0.07208'
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D10").Value = '''21.22'
$ws.Range("E10").Value = '  +2.42%  '
$ws.Range("D11").Value = '''0.9056'
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("E12").Value = '  +2.13%  '
$ws.Range("D13").Value = '1.885.46'
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("D14").Value = '''94.87'
$ws.Range("E14").Value = '  +2.62%  '
$ws.Range("D15").Value = '''5.268'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").Value = '''0.000008474'
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = '''14.45'
$ws.Range("E18").Value = '  +2.29%  '
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").Value = '27.100.02'
$ws.Range("E20").Value = '  +0.70%  '
$ws.Range("D21").Value = '''5.073'
$ws.Range("E21").Value = '  +0.91%  '
$ws.Range("D22").Value = '2.121.09'
$ws.Range("E22").Value = '  +1.50%  '
$ws.Range("E23").Value = '  +2.07%  '
$ws.Range("D24").Value = '''6.401'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = '''2.285'
$ws.Range("E25").Value = '  +10.96%  '
$ws.Range("D26").Value = '''146.00'
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("D27").Value = '''1.768'
$ws.Range("E27").Value = '  -1.38%  '
$ws.Range("D28").Value = '''18.06'
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("D29").Value = '''114.43'
$ws.Range("E29").Value = '  +1.13%  '
$ws.Range("D30").Value = '''4.947'
$ws.Range("E30").Value = '  +5.80%  '
$ws.Range("D31").Value = '''4.831'
$ws.Range("E31").Value = '  +3.18%  '
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").Value = '''0.05088'
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("D34").Value = '''1.236'
$ws.Range("E34").Value = '  +7.57%  '
$ws.Range("D35").Value = '''0.7815'
$ws.Range("E35").Value = '  +4.91%  '
$ws.Range("D36").Value = '''2.983'
$ws.Range("E36").Value = '  +1.21%  '
$ws.Range("D37").Value = '''3.285'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +5.08%  '
$ws.Range("D39").Value = '''0.01997'
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").Value = '''0.5591'
$ws.Range("E40").Value = '  +0.94%  '
$ws.Range("D41").Value = '''1.075'
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("D42").Value = '''9.108'
$ws.Range("E42").Value = '  +7.40%  '
$ws.Range("D43").Value = '''6.634'
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = '''0.1508'
$ws.Range("E45").Value = '  +2.62%  '
$ws.Range("D46").Value = '''0.4804'
$ws.Range("E46").Value = '  +2.71%  '
$ws.Range("D47").Value = '''10.22'
$ws.Range("E47").Value = '  +1.93%  '
$ws.Range("D48").Value = '''0.9987'
$ws.Range("E48").Value = '  -0.12%  '
$ws.Range("D49").Value = '''1.602'
$ws.Range("E49").Value = '  +2.46%  '
$ws.Range("D50").Value = '''37.58'
$ws.Range("E50").Value = '  +1.58%  '
$ws.Range("D51").Value = '''64.04'
$ws.Range("E51").Value = '  +1.57%  '
